# Ccl12-Ackr4.xlsx update: refresh NATMI edge-weight stats with the
# re-run ("new tpm") values, and drop the no-longer-produced MuSCs ->
# sending-cluster rows (the model now only reports FAPs and
# Resolving-Mac as sending clusters; MuSCs / Resolving-Mac remain as
# valid target clusters).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to have 9 data rows (r2:r10); it now only has 8
# (r2:r9). Drop the trailing row first so the used range / dimension
# shrinks to A1:T9, then overwrite the remaining data rows in place.
$ws.Rows(10).Delete()

$data = @(
    [PSCustomObject]@{ Row=2; A="FAPs"; B="Ccl12"; C="Ackr4"; D="ECs"; E=3; F=1; G=0.857148; H=2.571444; I=0.02394963654761903; J=0.02394963654761903; K=1; L=0.3333333333333333; M=0.1143813333333333; N=0.343144; O=0.128300337591142; P=0.1283003375911419; Q=0.09804173110400001; R=0.882375579936; S=0.003072746454244673; T=0.003072746454244672 },
    [PSCustomObject]@{ Row=3; A="FAPs"; B="Ccl12"; C="Ackr4"; D="FAPs"; E=3; F=1; G=0.857148; H=2.571444; I=0.02394963654761903; J=0.02394963654761903; K=3; L=1; M=0.7200953333333334; N=2.160286; O=0.8077233554817153; P=0.8077233554817151; Q=0.6172282747760001; R=5.555054472984001; S=0.01934468079481036; T=0.01934468079481036 },
    [PSCustomObject]@{ Row=4; A="FAPs"; B="Ccl12"; C="Ackr4"; D="MuSCs"; E=3; F=1; G=0.857148; H=2.571444; I=0.02394963654761903; J=0.02394963654761903; K=1; L=0.3333333333333333; M=0.03357866666666667; N=0.100736; O=0.03766483694187069; P=0.03766483694187069; Q=0.028781886976; R=0.259036982784; S=0.0009020591553831374; T=0.0009020591553831374 },
    [PSCustomObject]@{ Row=5; A="FAPs"; B="Ccl12"; C="Ackr4"; D="Resolving-Mac"; E=3; F=1; G=0.857148; H=2.571444; I=0.02394963654761903; J=0.02394963654761903; K=1; L=0.3333333333333333; M=0.023457; N=0.070371; O=0.02631146998527222; P=0.02631146998527222; Q=0.020106120636; R=0.180955085724; S=0.0006301501431808566; T=0.0006301501431808565 },
    [PSCustomObject]@{ Row=6; A="Resolving-Mac"; B="Ccl12"; C="Ackr4"; D="ECs"; E=3; F=1; G=34.93245566666667; H=104.797367; I=0.9760503634523809; J=0.9760503634523809; K=1; L=0.3333333333333333; M=0.1143813333333333; N=0.343144; O=0.128300337591142; P=0.1283003375911419; Q=3.995620855760889; R=35.96058770184801; S=0.1252275911368973; T=0.1252275911368972 },
    [PSCustomObject]@{ Row=7; A="Resolving-Mac"; B="Ccl12"; C="Ackr4"; D="FAPs"; E=3; F=1; G=34.93245566666667; H=104.797367; I=0.9760503634523809; J=0.9760503634523809; K=3; L=1; M=0.7200953333333334; N=2.160286; O=0.8077233554817153; P=0.8077233554817151; Q=25.15469830744023; R=226.392284766962; S=0.7883786746869048; T=0.7883786746869048 },
    [PSCustomObject]@{ Row=8; A="Resolving-Mac"; B="Ccl12"; C="Ackr4"; D="MuSCs"; E=3; F=1; G=34.93245566666667; H=104.797367; I=0.9760503634523809; J=0.9760503634523809; K=1; L=0.3333333333333333; M=0.03357866666666667; N=0.100736; O=0.03766483694187069; P=0.03766483694187069; Q=1.172985284679111; R=10.556867562112; S=0.03676277778648755; T=0.03676277778648755 },
    [PSCustomObject]@{ Row=9; A="Resolving-Mac"; B="Ccl12"; C="Ackr4"; D="Resolving-Mac"; E=3; F=1; G=34.93245566666667; H=104.797367; I=0.9760503634523809; J=0.9760503634523809; K=1; L=0.3333333333333333; M=0.023457; N=0.070371; O=0.02631146998527222; P=0.02631146998527222; Q=0.8194106125730002; R=7.374695513157001; S=0.02568131984209136; T=0.02568131984209136 }
)

foreach ($row in $data) {
    $r = $row.Row
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
    $ws.Range("K$r").Value = $row.K
    $ws.Range("L$r").Value = $row.L
    $ws.Range("M$r").Value = $row.M
    $ws.Range("N$r").Value = $row.N
    $ws.Range("O$r").Value = $row.O
    $ws.Range("P$r").Value = $row.P
    $ws.Range("Q$r").Value = $row.Q
    $ws.Range("R$r").Value = $row.R
    $ws.Range("S$r").Value = $row.S
    $ws.Range("T$r").Value = $row.T
}
